# Rearrange processing order and cleanup
# Updates the resident-passenger zonal distribution figures on both
# "Municipal Zone Summary" and "PMSA Summary" sheets to reflect the
# re-run / reprocessed numbers, and drops the stray decimal "weighted"
# figures in favor of clean re-derived integers.

$wb = $excel.ActiveWorkbook

$zoneData = @{
  2  = @(78, 155, 153, 115, 229, 227)
  3  = @(110, 214, 208, 201, 391, 380)
  4  = @(20, 40, 40, 33, 66, 66)
  5  = @(14, 27, 26, 19, 37, 36)
  6  = @(40, 77, 77, 71, 137, 137)
  7  = @(45, 89, 85, 63, 124, 118)
  8  = @(46, 92, 90, 83, 166, 164)
  9  = @(8, 16, 16, 10, 20, 20)
  10 = @(34, 66, 66, 46, 87, 87)
  11 = @(7, 14, 14, 11, 22, 22)
  12 = @(17, 33, 33, 35, 69, 69)
  13 = @(69, 135, 135, 95, 186, 186)
  14 = @(27, 53, 51, 42, 83, 77)
  15 = @(140, 274, 269, 243, 472, 463)
  16 = @(989, 1906, 1836, 1570, 3016, 2899)
  17 = @(29, 58, 58, 52, 104, 104)
  18 = @(10, 19, 18, 14, 27, 26)
  19 = @(6, 12, 12, 9, 18, 18)
  20 = @(13, 26, 26, 28, 56, 56)
  21 = @(166, 293, 254, 327, 591, 529)
}

$pmsaData = @{
  2  = @(81, 155, 146, 135, 257, 241)
  3  = @(442, 846, 805, 711, 1361, 1284)
  4  = @(528, 1028, 1006, 838, 1629, 1600)
  5  = @(141, 274, 267, 243, 469, 456)
  6  = @(164, 321, 317, 267, 522, 516)
  7  = @(232, 457, 450, 334, 658, 649)
  8  = @(110, 217, 214, 206, 402, 397)
  9  = @(4, 8, 8, 6, 12, 12)
  10 = @(166, 293, 254, 327, 591, 529)
}

$cols = @("B", "C", "D", "E", "F", "G")

$wsZone = $wb.Worksheets.Item("Municipal Zone Summary")
foreach ($row in $zoneData.Keys) {
  $values = $zoneData[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsZone.Range("$($cols[$i])$row").Value = $values[$i]
  }
}

$wsPmsa = $wb.Worksheets.Item("PMSA Summary")
foreach ($row in $pmsaData.Keys) {
  $values = $pmsaData[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsPmsa.Range("$($cols[$i])$row").Value = $values[$i]
  }
}
